$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the category labels in column A (per commit: "updated all the labels for review")
$ws.Range("A3").Value = "Insomnia related symptoms"
$ws.Range("A5").Value = "Insomnia related impairments"
$ws.Range("A6").Value = "Sleep behaviors"
